# #35 Updated todo history test case documentation with new test case to
# filter todo history by user.
#
# Adds a new test-case row (row 5) to the "Test Cases - ToDoHistory" sheet:
#   3.3 | Test displaying of individual to-do history for each user |
#   To test if the to-do history shown are the ones created by the user
#   and not other users. | - | After user log in, he/she sees only
#   his/her own to-do history. Logging into another user's account will
#   show a different list of to-do history. | - | Fail | -
#
# Also moves the active selection/viewport to E5 (matching the author's
# last on-screen selection when the row was filled in).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 was a blank, but already-bordered/formatted, placeholder row
# (same look as rows 2-4). Copy the formatting from the fully populated
# row 4 down onto row 5 first, so the new row picks up the same borders /
# alignment / wrap-text / number format as the rest of the table.
$ws.Range("A4:H4").Copy()
$ws.Range("A5:H5").PasteSpecial(-4122)

# Fill in the new test case's data.
$ws.Range("A5").Value = 3.3
$ws.Range("B5").Value = "Test displaying of individual to-do history for each user"
$ws.Range("C5").Value = "To test if the to-do history shown are the ones created by the user and not other users."
$ws.Range("D5").Value = "-"
$ws.Range("E5").Value = "After user log in, he/she sees only his/her own to-do history. Logging into another user's account will show a different list of to-do history."
$ws.Range("F5").Value = "-"
$ws.Range("G5").Value = "Fail"
$ws.Range("H5").Value = "-"

# Match the saved view state: scrolled back to column A and the E5 cell
# selected (instead of the old D1-scrolled / H1:H4 selection).
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("E5").Select() | Out-Null
